# Applies:
#  1) "day" sheet: D760:D773 bsecode values converted from text to numbers
#  2) "week" sheet: append 19 new rows (395-413) of stock data, extending
#     the used range from A1:I394 to A1:I413

$wb = $excel.ActiveWorkbook

# --- 1) Fix bsecode column (D) on the "day" sheet: text -> number ---
$day = $wb.Worksheets.Item("day")

$dayFixes = @{
    760 = 500387
    761 = 500124
    762 = 532488
    763 = 500825
    764 = 500300
    765 = 500483
    766 = 532286
    767 = 500260
    768 = 500440
    769 = 531642
    770 = 500670
    771 = 500038
    772 = 500085
    773 = 540691
}

foreach ($r in $dayFixes.Keys) {
    $day.Cells.Item($r, 4).Value = $dayFixes[$r]
}

# --- 2) Append new rows to the "week" sheet ---
$week = $wb.Worksheets.Item("week")

$newRows = @(
    @(1,  "SIEMENS",    "Siemens Limited",                                   "500550", 1.11,  7803,    274616),
    @(2,  "ALKEM",      "Alkem Laboratories Limited",                        "539523", 1.82,  6132.75, 253440),
    @(3,  "INDIGO",     "Interglobe Aviation Limited",                       "539448", 0.84,  4663.05, 624915),
    @(4,  "JKCEMENT",   "Jk Cement Limited",                                 "532644", 0.26,  4302.7,  84948),
    @(5,  "ESCORTS",    "Escorts Limited",                                   "500495", 0.7,   3815.4,  194550),
    @(6,  "INDIAMART",  "Indiamart Intermesh Ltd",                           "542726", -0.82, 3013.15, 303184),
    @(7,  "MUTHOOTFIN", "Muthoot Finance Limited",                           "533398", 0.68,  1968.75, 565171),
    @(8,  "AUROPHARMA", "Aurobindo Pharma Limited",                          "524804", 0.13,  1470.5,  437705),
    @(9,  "GODREJCP",   "Godrej Consumer Products Limited",                  "532424", -0.83, 1336.15, 658203),
    @(10, "IRCTC",      "Indian Railway Catering & Tourism Corporation Ltd", "542830", 1.06,  881,     892215),
    @(11, "CANFINHOME", "Can Fin Homes Limited",                             "511196", 0.28,  860,     175385),
    @(12, "SBIN",       "State Bank Of India",                               "500112", 1.15,  820.4,   11704698),
    @(13, "SBICARD",    "SBI Cards & Payment Services Ltd",                  "543066", 0.02,  740.15,  824220),
    @(14, "LAURUSLABS", "Laurus Labs Limited",                               "540222", 1.92,  475.15,  3771206),
    @(15, "NTPC",       "Ntpc Limited",                                      "532555", 1.72,  424.95,  8828971),
    @(16, "ABFRL",      "Aditya Birla Fashion And Retail Limited",           "535755", 1.01,  334.05,  1398893),
    @(17, "NMDC",       "Nmdc Limited",                                      "526371", 3.31,  231.44,  10451938),
    @(18, "FEDERALBNK", "The Federal Bank  Limited",                         "500469", 0.87,  195.3,   7308027),
    @(19, "CANBK",      "Canara Bank",                                       "532483", 2.12,  104.67,  19968661)
)

$startRow = 395
$dateTime = "18/10/2024 11:34:24"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $week.Cells.Item($r, 1).Value = $row[0]       # sr
    $week.Cells.Item($r, 2).Value = $row[1]       # nsecode
    $week.Cells.Item($r, 3).Value = $row[2]       # name

    # bsecode (column D) is kept as plain text in the source data (it is
    # never used numerically), so format it as Text before writing it,
    # then restore the "Normal" cell style so no stray number-format /
    # quote-prefix styling is left behind on the cell.
    $week.Cells.Item($r, 4).NumberFormat = "@"
    $week.Cells.Item($r, 4).Value = $row[3]
    $week.Cells.Item($r, 4).Style = "Normal"

    $week.Cells.Item($r, 5).Value = $row[4]       # per_chg
    $week.Cells.Item($r, 6).Value = $row[5]       # close
    $week.Cells.Item($r, 7).Value = $row[6]       # volume
    $week.Cells.Item($r, 8).Value = "week"        # timeframe
    $week.Cells.Item($r, 9).Value = $dateTime     # Date Time
}
